# Atualiza datasets e ajustes das ligas
# Rebuilds the "Nome do Time / ID do Time / Link do Time" table with four
# new teams inserted (alphabetically, by position) and all hyperlinks
# regenerated to match the new row order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final data, in the exact target row order (row 2 .. row 19).
$names = @(
    "bugredasmissões",
    "C R Juvenal",
    "Doug Leal F.C",
    "Esquadrão Gazembrino",
    "FBC Colorado",
    "GaúchoDaFronteira F.C",
    "GE Bebum",
    "GrioTeam",
    "Grêmio_Campeão_LA_27",
    "JV5 Tricolor Gaúcho",
    "La Primeira Patada Es Nuestra",
    "lsauer fc",
    "Medonho´s F.C.",
    "NHU PORÃ SAF.",
    "Pontaç0 F.C.",
    "SC 100 Sono",
    "SC ÉoINTER!",
    "Texas Club 2026"
)

$ids = @(
    19209079,
    1488983,
    287965,
    2916559,
    186283,
    2371918,
    16411206,
    14933455,
    47775950,
    1747619,
    32966,
    44810918,
    1867254,
    4088673,
    20651178,
    14709358,
    184499,
    1273719
)

# Wipe any hyperlinks on the sheet first -- they will be rebuilt from
# scratch below, in the final row order, so the relationship ids line up
# sequentially (rId1, rId2, ...) exactly like a freshly authored sheet.
$ws.Hyperlinks.Delete()

$firstDataRow = 2
$lastDataRow = $firstDataRow + $names.Length - 1

for ($i = 0; $i -lt $names.Length; $i++) {
    $row = $firstDataRow + $i
    $name = $names[$i]
    $id = $ids[$i]
    $link = "https://cartola.globo.com/#!/time/$id"

    $ws.Cells.Item($row, 1).Value = $name
    $ws.Cells.Item($row, 2).Value = $id
    $ws.Cells.Item($row, 3).Value = $link
}

# Remove any leftover rows below the new data (not needed here since the
# new table is longer than the old one, but keep the sheet tidy just in
# case a shorter edit is ever re-applied on top of a larger one).
$oldUsedRows = $ws.UsedRange.Rows.Count
if ($oldUsedRows -gt $lastDataRow) {
    $extra = $ws.Range($ws.Cells.Item($lastDataRow + 1, 1), $ws.Cells.Item($oldUsedRows, 3))
    $extra.ClearContents()
}

# Re-create the hyperlinks, one per data row, in order -- this reproduces
# rId1..rIdN sequential relationship ids, matching the original authoring
# pattern (Address = site root, SubAddress = the "!/time/<id>" fragment).
$linkStyle = $ws.Cells.Item(2, 3).Style
for ($i = 0; $i -lt $names.Length; $i++) {
    $row = $firstDataRow + $i
    $id = $ids[$i]
    $cell = $ws.Cells.Item($row, 3)
    $ws.Hyperlinks.Add($cell, "https://cartola.globo.com/", "!/time/$id")
    # Hyperlinks.Add() re-applies its own "Hyperlink" xf record; restore the
    # sheet's original hyperlink cell style so every C-column cell keeps
    # using the same style index the workbook already defined.
    $cell.Style = $linkStyle
}
